# Update the Android/Java study schedule (Sheet1):
#  - C14: fill in the completion date for the MaterialDesign/shell task
#  - Row 15: new task entry (task, expected date, and notes/解答)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- C14: completion date ---------------------------------------------
# Copy the date-formatted style from the neighbouring B14 cell first, so
# the new value inherits the same number format / borders as the rest of
# the "完成日期" column, then write the value.
$ws.Range("B14").Copy()
$ws.Range("C14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C14").Value = 43236

# --- Row 15: new task row ----------------------------------------------
# Row 15 was blank; pull the per-column formatting down from row 14 (the
# preceding populated row) before filling in the new task data.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A15").Value = "1.Android多媒体(Notification,MediaPlayer,VideoView)`n2.温习shell脚本(全面学习)`n"
$ws.Range("B15").Value = 43241
$ws.Range("D15").Value = "1.第一行代码第8章`n2.Linux命令行与shell脚本编程大全(第11章)，这里面要是有看不懂的要及时跟我说，或者随时提问题到word中，不用每次等到结束时才提。这本书是比较全面的，怕你看着难受。"

# --- Selection cursor ----------------------------------------------------
$ws.Range("E14").Select() | Out-Null
